# The "domesticParent.ueiDUNS" column (AG) is being removed from the CSV
# response example sheet, per "Updated changes for removal duns".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting the whole column shifts every subsequent column one to the left
# and removes the now-unused "domesticParent.ueiDUNS" shared string.
$ws.Columns("AG:AG").Delete()

# Reflect the resulting selection/view state (the column that now sits where
# AG used to be is selected, matching the post-edit cursor position).
$ws.Columns("AG:AG").Select()
